$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = $false
$ws.Range("D2").Value = 0.9999917431614516
$ws.Range("E2").Value = 0.9999917431614516

# Row 3
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 1

# Row 4
$ws.Range("C4").Value = $true
$ws.Range("D4").Value = 0.006985664340243368
$ws.Range("E4").Value = 0.006985664340243368

# Row 5
$ws.Range("D5").Value = [double]"8.410612112418685E-06"
$ws.Range("E5").Value = [double]"8.410612112418685E-06"

# Row 6
$ws.Range("D6").Value = [double]"1.07109936256667E-05"
$ws.Range("E6").Value = [double]"1.07109936256667E-05"

# Row 8
$ws.Range("D8").Value = 0.9998826443948218
$ws.Range("E8").Value = 0.0001173556051782265

# Row 10
$ws.Range("D10").Value = [double]"8.306214752864594E-16"
$ws.Range("E10").Value = 0.9999999999999992

# Row 11
$ws.Range("C11").Value = $false
$ws.Range("D11").Value = [double]"6.274378813596615E-08"
$ws.Range("E11").Value = 0.9999999372562118
$ws.Range("F11").Value = 230.6071472167969
$ws.Range("G11").Value = 0.5
